$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.776176691055298
$ws.Range("B1").Value = 1.766152143478394
$ws.Range("C1").Value = 2.087800264358521
$ws.Range("D1").Value = 3.892643928527832
$ws.Range("E1").Value = 3.96714973449707
